# Backlog.xlsx edit:
#  - Row 29 (requisito "O software deverá ter Wi-fi.") text changed to
#    "O software deverá ter uma rede de dados." (Luan's new requirement).
#  - Rows 25 & 26 reclassified from "Não Funcional" (red) to "Funcional"
#    (light green), matching the style already used by the other
#    "Funcional" rows (e.g. E12).
#  - Selection moved to B31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Replace the old "Wi-fi" requirement text with the new one.
$ws.Range("C29").Value = "O software deverá ter uma rede de dados."

# Reclassify rows 25 and 26 as "Funcional" (text + fill color, matching
# the look of the other Funcional cells such as E12).
$ws.Range("E25").Value = "Funcional"
$ws.Range("E25").Interior.Color = $ws.Range("E12").Interior.Color

$ws.Range("E26").Value = "Funcional"
$ws.Range("E26").Interior.Color = $ws.Range("E12").Interior.Color

# Move the active selection to B31.
[void]$ws.Range("B31").Select()
